$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.440.99"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.574.93"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.62"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3766"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.56%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.93"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3428"
$ws.Range("D9").ClearFormats()
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07682"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.30"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.009"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.943"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.575.97"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001137"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06759"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.84%  "
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.82"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.239"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.06"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.433"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "22.428.32"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.753"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -5.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.36"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "146.37"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.039"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.42"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.749.38"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.224"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("B33").Value = "WEMIXTOKEN"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.015"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.004"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.10"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08581"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02556"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2320"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06579"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.331"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +7.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.460"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.98%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.60"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.42%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6461"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.29%  "
$ws.Range("E44").Value = "  -1.98%  "
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.802"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6029"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.298"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +9.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.089"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "125.64"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07328"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.86%  "

Write-Host "Applied 103 cell updates"